$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.143.29"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.854.81"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.08"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.02"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0988"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "2.125.05"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.863.03"
$ws.Range("E13").Value = "  +4.34%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.44"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D17").Value = "35.114.65"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.88"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.80"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.21"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.66"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  +23.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.63"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0555"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +24.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.848"
$ws.Range("E35").Value = "  +20.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.02"
$ws.Range("E36").Value = "  +10.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.31"
$ws.Range("E37").Value = "  +6.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("E38").Value = "  +5.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.24"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "1.341.48"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.79"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.51"
$ws.Range("E45").Value = "  +45.88%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0555"
$ws.Range("E46").Value = "  +6.33%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.47"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").Value = "2.036.85"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  +0.21%  "
